# JS-Frameworks-Self-Evaluation-Protocol.xlsx
#
# Commit: "Delete Post Comment, Edit Post ADDED"
#
# The author gave a Score of 5 (out of the Maximum of 5 already in column D)
# to the two previously-ungraded checklist items:
#   - Row 28 "Edit Post"       -> C28 = 5
#   - Row 41 "Delete Comment"  -> C41 = 5
# The worksheet's Total Score (C44 = SUM(C6:C43)) recalculates from 231 to 241
# automatically as a consequence.
#
# (The diff also shows the author's window/scroll position moved while they
# were working - xWindow 2220->3330 and the visible top-left cell moving to
# A18 - these are just cosmetic view-state artifacts of the editing session,
# reproduced here on a best-effort basis.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self-Evaluation-Protocol")

# --- Data edits ---

# Row 28: "Edit Post" score -> 5
$ws.Range("C28").Value = 5

# Row 41: "Delete Comment" score -> 5
$ws.Range("C41").Value = 5

# Make sure the cached formula result for the Total Score (C44) is refreshed.
$excel.CalculateFull()

# --- Cosmetic view-state (best effort) ---
$ws.Activate() | Out-Null
$ws.Range("C41").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$excel.Left = 3330
